$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ck2.txt" row (row 3) is being removed entirely; the "ck3.txt" row
# (formerly row 4) shifts up to become row 3, carrying its existing
# formatting (style index 1) with it.
$ws.Rows("3").Delete()

# The "ck1.txt" row's "# Vin Values" (C2) is updated from 39 to 70.
$ws.Range("C2").Value = 70

# A new, empty, underlined-font placeholder cell is added at D5.
$ws.Range("D5").Font.Underline = $true

# Leave the selection on the newly added cell, matching the saved view.
$ws.Range("D5").Select()
